$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.924.24"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "3.484.56"
$ws.Range("E3").Value = "  -0.70%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.33"
$ws.Range("E5").Value = "  -1.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.49"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("E7").Value = "  +3.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "3.485.81"
$ws.Range("E9").Value = "  -0.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  -0.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.01"
$ws.Range("E11").Value = "  -1.65%  "

$ws.Range("E12").Value = "  -0.88%  "

$ws.Range("D13").Value = "4.085.73"
$ws.Range("E13").Value = "  -0.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.27"
$ws.Range("E14").Value = "  -0.87%  "

$ws.Range("E15").Value = "  -0.88%  "

$ws.Range("D16").Value = "67.963.90"
$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("E17").Value = "  -2.67%  "

$ws.Range("D18").Value = "3.486.83"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.21"
$ws.Range("E19").Value = "  -2.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.12"
$ws.Range("E20").Value = "  -4.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "395.38"
$ws.Range("E21").Value = "  -0.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.94"
$ws.Range("E22").Value = "  -1.71%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.84"
$ws.Range("E23").Value = "  +2.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.540"
$ws.Range("E24").Value = "  -0.84%  "

$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "72.18"
$ws.Range("E26").Value = "  -2.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000122"
$ws.Range("E27").Value = "  -1.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.45"
$ws.Range("E28").Value = "  +0.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.177"
$ws.Range("E29").Value = "  -1.39%  "

$ws.Range("E30").Value = "  +3.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  -2.64%  "

$ws.Range("E32").Value = "  -1.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.40"
$ws.Range("E33").Value = "  -3.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.62"
$ws.Range("E34").Value = "  -1.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.37"
$ws.Range("E35").Value = "  -0.67%  "

$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("E37").Value = "  -6.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.72"
$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("E39").Value = "  +1.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  +5.00%  "

$ws.Range("E41").Value = "  -4.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.78"
$ws.Range("E42").Value = "  -4.80%  "

$ws.Range("E43").Value = "  -1.46%  "

$ws.Range("E44").Value = "  -1.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0722"
$ws.Range("E45").Value = "  -1.92%  "

$ws.Range("D46").Value = "2.749.97"
$ws.Range("E46").Value = "  -2.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.30"
$ws.Range("E47").Value = "  -5.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.50"
$ws.Range("E48").Value = "  -2.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0299"
$ws.Range("E49").Value = "  -1.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "330.69"

$ws.Range("E51").Value = "  -3.88%  "
